$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.629057
$ws.Range("H2").Value = 34.887171
$ws.Range("I2").Value = 0.3062678464977661
$ws.Range("J2").Value = 0.3062678464977662
$ws.Range("M2").Value = 34.49888633333333
$ws.Range("N2").Value = 103.496659
$ws.Range("O2").Value = 0.4998067520528027
$ws.Range("P2").Value = 0.4998067520528027
$ws.Range("Q2").Value = 401.1895156068542
$ws.Range("R2").Value = 3610.705640461688
$ws.Range("S2").Value = 0.1530747376162548
$ws.Range("T2").Value = 0.1530747376162549

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.629057
$ws.Range("H3").Value = 34.887171
$ws.Range("I3").Value = 0.3062678464977661
$ws.Range("J3").Value = 0.3062678464977662
$ws.Range("N3").Value = 72.35583600000001
$ws.Range("O3").Value = 0.3494212830891987
$ws.Range("P3").Value = 0.3494212830891987
$ws.Range("Q3").Value = 280.4767137088839
$ws.Range("R3").Value = 2524.290423379956
$ws.Range("S3").Value = 0.1070165038922152
$ws.Range("T3").Value = 0.1070165038922152

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.629057
$ws.Range("H4").Value = 34.887171
$ws.Range("I4").Value = 0.3062678464977661
$ws.Range("J4").Value = 0.3062678464977662
$ws.Range("M4").Value = 10.406952
$ws.Range("N4").Value = 31.220856
$ws.Range("O4").Value = 0.1507719648579985
$ws.Range("P4").Value = 0.1507719648579985
$ws.Range("Q4").Value = 121.023038004264
$ws.Range("R4").Value = 1089.207342038376
$ws.Range("S4").Value = 0.04617660498929609
$ws.Range("T4").Value = 0.04617660498929609

$ws.Range("I5").Value = 0.6269156120645606
$ws.Range("J5").Value = 0.6269156120645607
$ws.Range("M5").Value = 34.49888633333333
$ws.Range("N5").Value = 103.496659
$ws.Range("O5").Value = 0.4998067520528027
$ws.Range("P5").Value = 0.4998067520528027
$ws.Range("Q5").Value = 821.2157221420568
$ws.Range("R5").Value = 7390.941499278511
$ws.Range("S5").Value = 0.3133366558771829
$ws.Range("T5").Value = 0.313336655877183

$ws.Range("I6").Value = 0.6269156120645606
$ws.Range("J6").Value = 0.6269156120645607
$ws.Range("N6").Value = 72.35583600000001
$ws.Range("O6").Value = 0.3494212830891987
$ws.Range("P6").Value = 0.3494212830891987
$ws.Range("Q6").Value = 574.1223985977388
$ws.Range("R6").Value = 5167.101587379649
$ws.Range("S6").Value = 0.2190576575562491
$ws.Range("T6").Value = 0.2190576575562491

$ws.Range("I7").Value = 0.6269156120645606
$ws.Range("J7").Value = 0.6269156120645607
$ws.Range("M7").Value = 10.406952
$ws.Range("N7").Value = 31.220856
$ws.Range("O7").Value = 0.1507719648579985
$ws.Range("P7").Value = 0.1507719648579985
$ws.Range("Q7").Value = 247.728361994112
$ws.Range("R7").Value = 2229.555257947008
$ws.Range("S7").Value = 0.09452129863112857
$ws.Range("T7").Value = 0.0945212986311286

$ws.Range("G8").Value = 2.537038666666667
$ws.Range("H8").Value = 7.611116
$ws.Range("I8").Value = 0.06681654143767324
$ws.Range("J8").Value = 0.06681654143767324
$ws.Range("M8").Value = 34.49888633333333
$ws.Range("N8").Value = 103.496659
$ws.Range("O8").Value = 0.4998067520528027
$ws.Range("P8").Value = 0.4998067520528027
$ws.Range("Q8").Value = 87.52500858460489
$ws.Range("R8").Value = 787.725077261444
$ws.Range("S8").Value = 0.03339535855936496
$ws.Range("T8").Value = 0.03339535855936496

$ws.Range("G9").Value = 2.537038666666667
$ws.Range("H9").Value = 7.611116
$ws.Range("I9").Value = 0.06681654143767324
$ws.Range("J9").Value = 0.06681654143767324
$ws.Range("N9").Value = 72.35583600000001
$ws.Range("O9").Value = 0.3494212830891987
$ws.Range("P9").Value = 0.3494212830891987
$ws.Range("Q9").Value = 61.18985123033067
$ws.Range("R9").Value = 550.708661072976
$ws.Range("S9").Value = 0.02334712164073439
$ws.Range("T9").Value = 0.0233471216407344

$ws.Range("G10").Value = 2.537038666666667
$ws.Range("H10").Value = 7.611116
$ws.Range("I10").Value = 0.06681654143767324
$ws.Range("J10").Value = 0.06681654143767324
$ws.Range("M10").Value = 10.406952
$ws.Range("N10").Value = 31.220856
$ws.Range("O10").Value = 0.1507719648579985
$ws.Range("P10").Value = 0.1507719648579985
$ws.Range("Q10").Value = 26.402839626144
$ws.Range("R10").Value = 237.625556635296
$ws.Range("S10").Value = 0.01007406123757387
$ws.Range("T10").Value = 0.01007406123757387
